$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 1.75
$ws.Range("H3").Value = 3.5
$ws.Range("I3").Value = 5
$ws.Range("J3").Value = 2.4
$ws.Range("U3").Value = 2
$ws.Range("V3").Value = 1.75
$ws.Range("AG3").Value = 11
$ws.Range("AH3").Value = 23
$ws.Range("AN3").Value = 3.6
$ws.Range("AO3").Value = 9.5

# Row 4
$ws.Range("G4").Value = 1.95
$ws.Range("H4").Value = 3
$ws.Range("J4").Value = 2.75
$ws.Range("U4").Value = 2.63
$ws.Range("V4").Value = 1.44
$ws.Range("AD4").Value = 6.5

# Row 5
$ws.Range("K5").Value = 1.91
$ws.Range("O5").Value = 1.53
$ws.Range("P5").Value = 2.38
$ws.Range("BA5").Value = 126
$ws.Range("BB5").Value = 401

# Row 10
$ws.Range("J10").Value = 2.87
